# Balochi "Batal" proverbs workbook — append 7 more proverb rows (117-123)
# to Sheet1, mirroring the existing A (batal/proverb) / B (meaning) /
# C (source) layout, then move the view/selection to the new bottom row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows' data: (batal, meaning). Column C always repeats the same
# "source" shared string that every other data row already uses, so we
# copy it forward rather than retyping it.
$newRows = @(
    @("پدگروک نام نہ کَٹیت۔", "بے زانت ءُ پشت کپوک بے نام بنت نام نہ کٹ انت"),
    @("پدگِروک، پداکپیت۔", "آ کہ لالچ ءُ ضدّکنت گڑا پَشت کپیت دیم ءَ نہ روت"),
    @("پد، پہ پد بین۔", "کسے ءِ پُشت درد پُشت کجام وڑا بو تگاں آباد یا بزگ، گڑا آہانی آؤکیں نسل ہما وڑ بنت دانکہ چنکس پُشت بہ رؤت"),
    @("پرشتگیں کمان کار نہ دنت۔", "آ کہ بے جوہر ءُ بے کار اِنت آکار‌‌ ءَ نئیت"),
    @("پڑوکیں جَن شہ جِنّ ءَ گندہ ترانت۔", "آ جَن کہ فسادی اِنت دغا اِنت چہ آجن ءَ بے جَنی گہتر اِنت"),
    @("پزّوری، شہ زوری۔", "جوان ءُ بلد بئے گڑا پزّور بئے، زوراک بئے"),
    @("پزُوُل ہرچ ءِ دولت و بدپہریز ءِ صحت نہ مانیت۔", "بے پروائی پکیر کنت، ہما کہ دپ نہ پہلیت نادرہ بیت")
)

$lastRow = 116
$sourceValue = $ws.Range("C$lastRow").Value2

$row = $lastRow
foreach ($pair in $newRows) {
    $row = $row + 1

    # Clone row 116's direct formatting (font/alignment/border/row-height)
    # onto the new row before writing values, so the new rows look just
    # like the rest of the table instead of picking up default formatting.
    $ws.Range("A$lastRow`:C$lastRow").Copy()
    $ws.Range("A$row`:C$row").PasteSpecial(-4122) | Out-Null
    $ws.Rows.Item($row).RowHeight = $ws.Rows.Item($lastRow).RowHeight

    $ws.Range("A$row").Value2 = $pair[0]
    $ws.Range("B$row").Value2 = $pair[1]
    $ws.Range("C$row").Value2 = $sourceValue
}

$excel.CutCopyMode = 0

# Update the window to scroll/select like the author left it: top-left
# anchored at B112, with B123 (the new last row's meaning cell) selected.
$excel.ActiveWindow.ScrollRow = 112
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B$row").Select()
